$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.973.46"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.29%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.921.96"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.19%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.90"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.05%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4583"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.14%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3812"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.16%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07747"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.10%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9780"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.26%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "22.54"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.25%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.930.54"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.29%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.700"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.53%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.966"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.11%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06992"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.05%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "84.74"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.80%  "

$ws.Range("E17").Value = "  -0.11%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009489"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.54%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.70"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.19%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.004"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.15%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "28.978.69"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.45%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.341"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.08%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.08"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.55%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.180.48"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.31%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.057"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.97%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "158.17"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.53%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.09"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.38%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.617"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.57%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "117.56"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.08%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.840"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.15%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09284"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.18%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8642"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.25%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.106"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.14%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.245"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.60%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.015"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.17%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05689"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.06%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.152"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.53%  "

$ws.Range("E38").Value = "  -0.10%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02051"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.82%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.100"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +13.92%  "

$ws.Range("E41").Value = "  -0.16%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5506"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.09%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1757"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.09%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "9.332"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.11%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.000002772"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +17.95%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.181"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.76%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5182"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.21%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06938"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.57%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "11.18"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.66%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "110.84"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.64%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.765"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.36%  "

